$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Item ID 5489)
$ws.Range("H2").Value = 919.6
$ws.Range("J2").Value = 899
$ws.Range("L2").Value = 899
$ws.Range("N2").Value = -1125

# Row 5 (Item ID 5503)
$ws.Range("H5").Value = 837.4286
$ws.Range("I5").Value = 643.6667
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 643.6667
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -528.6667
$ws.Range("N5").Value = -2230

# Row 17 (Item ID 38956)
$ws.Range("H17").Value = 3270
$ws.Range("I17").Value = 3133.3333
$ws.Range("K17").Value = 9399.999899999999
$ws.Range("M17").Value = -9231.999899999999

# Row 19 (Item ID 7015)
$ws.Range("H19").Value = 999
$ws.Range("I19").Value = 999
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 999
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -824
$ws.Range("N19").ClearContents()

# Row 40 (Item ID 5505)
$ws.Range("H40").Value = 2599.4546
$ws.Range("I40").Value = 3475
$ws.Range("J40").Value = 2099.1428
$ws.Range("K40").Value = 3475
$ws.Range("L40").Value = 2099.1428
$ws.Range("M40").Value = -3300
$ws.Range("N40").Value = -2449.1428

# Row 138 (Item ID 44169)
$ws.Range("H138").Value = 3172.6667
$ws.Range("J138").Value = 3399.4075
$ws.Range("L138").Value = 10198.2225
$ws.Range("N138").Value = -20478.2225

$ws = $wb.Worksheets.Item("ARM")
# Row 29 (Item ID 3557)
$ws.Range("H29").Value = 1399.6666
$ws.Range("I29").Value = 1399.6666
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1399.6666
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1091.6666
$ws.Range("N29").ClearContents()

# Row 74 (Item ID 44000)
$ws.Range("H74").Value = 3600

# Row 77 (Item ID 44000)
$ws.Range("H77").Value = 3600

# Row 109 (Item ID 25646)
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774

# Row 110 (Item ID 27708)
$ws.Range("H110").Value = 62501500
$ws.Range("I110").Value = 1999.6666
$ws.Range("K110").Value = 1999.6666
$ws.Range("M110").Value = 45.33339999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Item ID 19939)
$ws.Range("H94").Value = 111368.3
$ws.Range("I94").Value = 111368.3
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 111368.3
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -110917.3
$ws.Range("N94").ClearContents()

# Row 138 (Item ID 42308)
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Item ID 27691)
$ws.Range("H16").Value = 448.75
$ws.Range("J16").Value = 365.66666
$ws.Range("L16").Value = 365.66666
$ws.Range("N16").Value = -939.66666

# Row 58 (Item ID 44021)
$ws.Range("H58").Value = 1478.091
$ws.Range("I58").Value = 1473
$ws.Range("K58").Value = 1473
$ws.Range("M58").Value = -1270

# Row 70 (Item ID 12011)
$ws.Range("H70").Value = 30000
$ws.Range("I70").Value = 30000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685

# Row 73 (Item ID 12011)
$ws.Range("H73").Value = 30000
$ws.Range("I73").Value = 30000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908

# Row 106 (Item ID 18661)
$ws.Range("H106").Value = 76714.28999999999
$ws.Range("J106").Value = 76714.28999999999
$ws.Range("L106").Value = 76714.28999999999
$ws.Range("N106").Value = -79238.28999999999

# Row 113 (Item ID 27691)
$ws.Range("H113").Value = 448.75
$ws.Range("J113").Value = 365.66666
$ws.Range("L113").Value = 365.66666
$ws.Range("N113").Value = -4705.66666

# Row 122 (Item ID 36196)
$ws.Range("H122").Value = 1601.4286
$ws.Range("I122").Value = 1601.4286
$ws.Range("K122").Value = 4804.2858
$ws.Range("M122").Value = -2354.2858

# Row 134 (Item ID 44020)
$ws.Range("H134").Value = 2710.25
$ws.Range("J134").Value = 2425
$ws.Range("L134").Value = 7275
$ws.Range("N134").Value = -12345

# Row 136 (Item ID 44021)
$ws.Range("H136").Value = 1478.091
$ws.Range("I136").Value = 1473
$ws.Range("K136").Value = 4419
$ws.Range("M136").Value = -1869

$ws = $wb.Worksheets.Item("CUL")
# Row 36 (Item ID 4732)
$ws.Range("H36").Value = 676.6
$ws.Range("I36").Value = 676.6
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2029.8
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1860.8
$ws.Range("N36").ClearContents()

# Row 37 (Item ID 9516)
$ws.Range("H37").Value = 27500
$ws.Range("J37").Value = 27500
$ws.Range("L37").Value = 82500
$ws.Range("N37").Value = -82724

# Row 75 (Item ID 12863)
$ws.Range("H75").Value = 962.3333
$ws.Range("J75").Value = 1039.2222
$ws.Range("L75").Value = 3117.6666
$ws.Range("N75").Value = -5113.6666

# Row 78 (Item ID 12863)
$ws.Range("H78").Value = 962.3333
$ws.Range("J78").Value = 1039.2222
$ws.Range("L78").Value = 9352.9998
$ws.Range("N78").Value = -19336.9998

# Row 131 (Item ID 36060)
$ws.Range("H131").Value = 1647.5333
$ws.Range("I131").Value = 781.8
$ws.Range("K131").Value = 2345.4
$ws.Range("M131").Value = 2694.6

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Item ID 12521)
$ws.Range("H80").Value = 4047.7
$ws.Range("I80").Value = 2839.4
$ws.Range("J80").Value = 5256
$ws.Range("K80").Value = 2839.4
$ws.Range("L80").Value = 5256
$ws.Range("M80").Value = -1841.4
$ws.Range("N80").Value = -7252

# Row 83 (Item ID 12521)
$ws.Range("H83").Value = 4047.7
$ws.Range("I83").Value = 2839.4
$ws.Range("J83").Value = 5256
$ws.Range("K83").Value = 14197
$ws.Range("L83").Value = 26280
$ws.Range("M83").Value = -9205
$ws.Range("N83").Value = -36264

# Row 101 (Item ID 18513)
$ws.Range("H101").Value = 9850.666999999999
$ws.Range("J101").Value = 9850.666999999999
$ws.Range("L101").Value = 9850.666999999999
$ws.Range("N101").Value = -16340.667

# Row 132 (Item ID 44008)
$ws.Range("H132").Value = 2166.6667
$ws.Range("I132").Value = 2250
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Item ID 5289)
$ws.Range("H16").Value = 1362.5
$ws.Range("I16").Value = 1362.5
$ws.Range("K16").Value = 1362.5
$ws.Range("M16").Value = -1192.5

# Row 22 (Item ID 5277)
$ws.Range("H22").Value = 213
$ws.Range("I22").Value = 192
$ws.Range("K22").Value = 192
$ws.Range("M22").Value = 103

# Row 27 (Item ID 5277)
$ws.Range("H27").Value = 213
$ws.Range("I27").Value = 192
$ws.Range("K27").Value = 192
$ws.Range("M27").Value = -85

# Row 46 (Item ID 5282)
$ws.Range("H46").Value = 888
$ws.Range("J46").Value = 888
$ws.Range("L46").Value = 888
$ws.Range("N46").Value = -1264

# Row 51 (Item ID 3423)
$ws.Range("H51").Value = 42000
$ws.Range("I51").Value = 42000
$ws.Range("K51").Value = 42000
$ws.Range("M51").Value = -41522

# Row 93 (Item ID 19993)
$ws.Range("H93").Value = 25642362
$ws.Range("I93").Value = 33334724
$ws.Range("K93").Value = 33334724
$ws.Range("M93").Value = -33333476

# Row 94 (Item ID 18067)
$ws.Range("H94").Value = 43328.668
$ws.Range("J94").Value = 43328.668
$ws.Range("L94").Value = 43328.668
$ws.Range("N94").Value = -44680.668

# Row 136 (Item ID 44060)
$ws.Range("H136").Value = 2849.5
$ws.Range("J136").Value = 2642.1428
$ws.Range("L136").Value = 7926.428400000001
$ws.Range("N136").Value = -13026.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 22 (Item ID 3041)
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 24 (Item ID 3561)
$ws.Range("H24").Value = 19997.5
$ws.Range("I24").Value = 19997.5
$ws.Range("K24").Value = 19997.5
$ws.Range("M24").Value = -19767.5

# Row 122 (Item ID 36208)
$ws.Range("H122").Value = 1901.3334
$ws.Range("I122").Value = 1004
$ws.Range("K122").Value = 3012
$ws.Range("M122").Value = -562

# Row 123 (Item ID 34127)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
